$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# New task rows (10-15): task name (A), role (B, where applicable),
# date (C), start time (D) and end time (E). Columns F-H already carry
# their formulas/formatting and recalc automatically once D/E are filled.

$ws.Range("A10").Value = "Udfør brugertest af Android mockup"
$ws.Range("C10").Value = 43887
$ws.Range("D10").Value = 0.4375
$ws.Range("E10").Value = 0.49305555555555558

$ws.Range("A11").Value = "Indsaml UC05 Beregn Markedsføringsbidrag"
$ws.Range("B11").Value = "Requirement Specifier"
$ws.Range("C11").Value = 43887
$ws.Range("D11").Value = 0.52083333333333337
$ws.Range("E11").Value = 0.5625

$ws.Range("A12").Value = "Lav UC05 Beregn Markedsføringsbidrag"
$ws.Range("B12").Value = "Requirement Specifier"
$ws.Range("C12").Value = 43887
$ws.Range("D12").Value = 0.5625
$ws.Range("E12").Value = 0.58333333333333337

$ws.Range("A13").Value = "Lav DOM05 Beregn Markedsføringsbidrag"
$ws.Range("B13").Value = "Requirement Specifier"
$ws.Range("C13").Value = 43887
$ws.Range("D13").Value = 0.58333333333333337
$ws.Range("E13").Value = 0.60416666666666663

$ws.Range("A14").Value = "Review move og rename vejledning"
$ws.Range("B14").Value = "Reviewer"
$ws.Range("C14").Value = 43887
$ws.Range("D14").Value = 0.60416666666666663
$ws.Range("E14").Value = 0.61458333333333337

$ws.Range("A15").Value = "Review AD02 Beregn vareforbrug"
$ws.Range("B15").Value = "Reviewer"
$ws.Range("C15").Value = 43887
$ws.Range("D15").Value = 0.61458333333333337
$ws.Range("E15").Value = 0.625

# Update the selected cell on the sheet (matches the saved view state).
$ws.Range("C20").Select()
